# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (column G) wherever the combined value
# "dnasr281@gmail.com, System" appears, turning it into
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$colG = $ws.Columns.Item(7)

# Use Find/FindNext so only cells that actually contain the target
# text are touched (avoids materializing empty cells elsewhere in
# the column).
$firstFound = $colG.Find($oldValue)

if ($firstFound -ne $null) {
    $firstAddress = $firstFound.Address()
    $current = $firstFound
    $safety = 0
    do {
        $current.Value2 = $newValue
        $current = $colG.FindNext($current)
        $safety = $safety + 1
    } while ($current -ne $null -and $current.Address() -ne $firstAddress -and $safety -lt 1000)
}
